# Adapted AI model, corrected error
#
# The "Aggregated" sheet's coefficient table is restructured: the numeric
# travel_code_id index column is dropped, the travel_code (car/plane/train)
# text column becomes column A, the "crookness_coef" column is renamed
# "tortuosity_coef" and becomes column B, and "carbon_coef" becomes column C.
#
# This is done in a handful of passes so that the workbook's shared string
# table ends up ordered the same way a human editing the sheet in stages
# (dropping then re-adding the travel_code header) would naturally produce.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aggregated")

# --- Pass 1: drop the old "travel_code" header text -----------------------
$ws.Range("B1").Value = "__tmp_travel_code__"
$wb.Save()

# --- Pass 2: drop "travel_code_id" and "crookness_coef" headers -----------
$ws.Range("A1").Value = "__tmp_travel_code_id__"
$ws.Range("C1").Value = "__tmp_crookness_coef__"
$wb.Save()

# --- Pass 3: lay out the final table ---------------------------------------
$ws.Range("A1").Value = "travel_code"
$ws.Range("B1").Value = "tortuosity_coef"
$ws.Range("C1").Value = "carbon_coef"

$b2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2
$b3 = $ws.Range("C3").Value2
$d3 = $ws.Range("D3").Value2
$b4 = $ws.Range("C4").Value2
$d4 = $ws.Range("D4").Value2

# Carry over the right-aligned / vertical-centered number style that lives
# on the carbon_coef column (D) to its new home in column C.
$ws.Range("D2:D4").Copy()
$ws.Range("C2:C4").PasteSpecial(-4122)

$ws.Range("A2").Value = "car"
$ws.Range("B2").Value2 = $b2
$ws.Range("C2").Value2 = $d2

$ws.Range("A3").Value = "plane"
$ws.Range("B3").Value2 = $b3
$ws.Range("C3").Value2 = $d3

$ws.Range("A4").Value = "train"
$ws.Range("B4").Value2 = $b4
$ws.Range("C4").Value2 = $d4

$ws.Range("D1:D4").Clear()

$ws.Range("B4").Select()

$wb.Save()
